$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update District column (G) values from "Tumkur" to "Tumakuru (Tumkur)"
# Only exact matches of "Tumkur" are updated; other district-like text (e.g. "Tumkuru",
# "Kunigal", or full school-name strings that ended up in column G) are left untouched.
for ($r = 1; $r -le 63; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Tumkur") {
        $cell.Value2 = "Tumakuru (Tumkur)"
    }
}

# Clear out the stray empty inline-string cells in column F (rows 14, 21, 23, 38, 57, 59)
# so they no longer exist as cells in the sheet.
$emptyFRows = @(14, 21, 23, 38, 57, 59)
foreach ($r in $emptyFRows) {
    $ws.Cells.Item($r, 6).Clear()
}
